# Auto update stock data
#
# The sheet lists several stocks, each as a 6-row block. The first row of
# every block holds "today's" date in column A (the other 5 rows are fixed
# historical fiscal year-end dates). This refreshes that date from
# 2025/11/22 to 2025/11/23 for every block, leaving all other cells and
# formatting untouched.
#
# Column A is plain text (e.g. "2024/12/31" lives right next to it as an
# untouched row), so a naive `$cell.Value = "2025/11/23"` would make Excel's
# type-inference reinterpret the new text as a date serial instead of
# leaving it as a string. To avoid that, stage the literal text in a
# scratch cell formatted as Text, then copy/PasteSpecial just the values
# into each target cell (PasteSpecial Values doesn't bring the scratch
# cell's number format along, so the destination's existing formatting is
# preserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "2025/11/22"
$newDate = "2025/11/23"

$usedRows = $ws.UsedRange.Rows.Count
$lastRow = $ws.UsedRange.Row + $usedRows - 1

# Scratch cell well below the data, used only to stage the replacement text
# so it round-trips through Excel as a string rather than being parsed as a
# date.
$helper = $ws.Cells.Item($lastRow + 50, 1)
$helper.NumberFormat = "@"
$helper.Value = $newDate
$helper.Copy()

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldDate) {
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
$helper.Clear()
